$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new row 25: "224. Basic Calculator" entry under a new "Stack & Recursion" topic ---

$ws.Range("A25").Value = 'Stack & Recursion'
$ws.Range("B25").Value = '224. Basic Calculator'
$ws.Range("C25").Value = 'first convert to postfix, in postfix put "," before everytime u encounter a operator to indicate end of a number
then in soultution everytime u encounter a single "," parse the final number u made from chars to Integer and store in stack
if u encounter 2 ",," simultaneously that means there was a a expression like 1 -(-2) which means u should push a 0 to the stack when u encounter a ",,"
if u encounter ,2,3+ in postfix this means the expression was "-2+3", in postfix we add "," evertime we encounter a operator to indicate end of num, since there was no num before "-" this means that the "-" is unary operator, so push a 0 here too!'
$ws.Range("D25").Value = 'My approach is very inefficent, checkout a efficient approach on YT or Leetcode'

# Row height to match the wrapped multi-line content (same as other "big text" rows, e.g. row 15)
$ws.Rows.Item(25).RowHeight = 115.2

# Formatting common to the whole new row: thin border all round, left/center aligned, wrap text
# (matches the formatting convention used by every other data row in the sheet)
$rowRange = $ws.Range("A25:D25")
$rowRange.Borders.LineStyle = 1
$rowRange.HorizontalAlignment = -4131
$rowRange.VerticalAlignment = -4108
$rowRange.WrapText = $true

# B25 ("224. Basic Calculator") is flagged in red, same convention as other flagged
# questions in column B (e.g. B22 - "84. Largest Rectangle in Histogram")
$ws.Range("B25").Interior.Color = 255

# D25 is a bold "note to self", same convention as the bold header row
$ws.Range("D25").Font.Bold = $true

# --- Update selection to reflect where the user ended up after the edit ---
$ws.Range("B27").Select()
